# Applies the weekly cryptos-list price/volume refresh described by the
# commit "Updated cryptos list on Thu Jun  1 23:50:48 UTC 2023 with GitHub
# Actions" -- per-cell literal value updates on Sheet1 (rows 2-51).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.837.01"
$ws.Range("E2").Value = "  -1.49%  "
$ws.Range("D3").Value = "1.861.87"
$ws.Range("E3").Value = "  -0.66%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "305.18"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.56%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9999"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5079"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.10%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3652"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.41%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07153"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.14%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8885"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.68%  "
$ws.Range("E11").Value = "  -1.43%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07476"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.78%  "
$ws.Range("D13").Value = "1.866.60"
$ws.Range("E13").Value = "  -0.49%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "93.99"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.68%  "
$ws.Range("E15").Value = "  -1.96%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.000"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.14%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008482"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.49%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.14"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.11%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.000"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.06%  "
$ws.Range("D20").Value = "26.876.33"
$ws.Range("E20").Value = "  -1.43%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.988"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.53%  "
$ws.Range("D22").Value = "2.111.23"
$ws.Range("E22").Value = "  +0.02%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.33"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.32%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.346"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.13%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.78"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.71%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.773"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.40%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.85"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.90%  "
$ws.Range("E28").Value = "  -0.27%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "113.37"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.05%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.668"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.07%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.694"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.04%  "
$ws.Range("E32").Value = "  -1.56%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05025"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.21%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7469"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.44%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.931"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.05%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.149"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.32%  "
$ws.Range("E37").Value = "  +2.11%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.498"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.22%  "
$ws.Range("B39").Value = "TheSandbox"
$ws.Range("C39").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5548"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.19%  "
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01976"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.93%  "
$ws.Range("E41").Value = "  -0.32%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.559"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.18%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "115.86"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.90%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.578"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.55%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1482"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.27%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4734"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.95%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.9999"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.06%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.946"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.81%  "
$ws.Range("E49").Value = "  +0.50%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.550"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.13%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "62.74"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.79%  "
